# testData.xlsx / Sheet1: extend the user list.
# Row 2's username "user0" becomes "user1"; three more rows are appended
# below it (user2, user3, user4), each repeating the same password/demo
# columns as row 2 (John / demo) - the classic "type a value, then fill
# down" pattern.  Finish with C2:C5 selected (mirrors the saved selection
# in the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "user1"
$ws.Range("B2").Value = "john"
$ws.Range("C2").Value = "demo"

$ws.Range("A3").Value = "user2"
$ws.Range("B3").Value = "john"
$ws.Range("C3").Value = "demo"

$ws.Range("A4").Value = "user3"
$ws.Range("B4").Value = "john"
$ws.Range("C4").Value = "demo"

$ws.Range("A5").Value = "user4"
$ws.Range("B5").Value = "john"
$ws.Range("C5").Value = "demo"

$ws.Range("C2:C5").Select()
